$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemoWebShop")

$ws.Range("C2").Value = "yYWfVclu@gmail.com"
$ws.Range("F2").Value = "GkVWM"
$ws.Range("G2").Value = "dIqyl"
